# Daily attendance processing - 2025-12-05 11:48:35
# Normalize the "Recorded By" (column G) entries so that "System" (capitalized)
# is listed first among the comma-separated recorder names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "system, System, backup@backdoor.com") {
        $cell.Value = "System, system, backup@backdoor.com"
    }
}
